$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text (or $null to leave unchanged), new Volume(1h) (E) text (or $null to leave unchanged).
$updates = @(
    @{ Row = 2; D = "67.101.90"; E = "  -0.35%  " }
    @{ Row = 3; D = "2.613.49"; E = "  -0.90%  " }
    @{ Row = 4; D = $null; E = "  -0.04%  " }
    @{ Row = 5; D = "591.22"; E = "  -1.07%  " }
    @{ Row = 6; D = "166.13"; E = "  -0.03%  " }
    @{ Row = 7; D = $null; E = "  +0.05%  " }
    @{ Row = 8; D = "0.532"; E = "  -1.89%  " }
    @{ Row = 9; D = "2.613.41"; E = "  -0.91%  " }
    @{ Row = 10; D = "0.138"; E = "  -4.51%  " }
    @{ Row = 12; D = "0.363"; E = "  -0.13%  " }
    @{ Row = 13; D = "5.20"; E = "  -0.22%  " }
    @{ Row = 14; D = "27.33"; E = "  -2.03%  " }
    @{ Row = 15; D = "3.089.13"; E = "  -1.00%  " }
    @{ Row = 16; D = $null; E = "  -2.18%  " }
    @{ Row = 17; D = "67.290.98"; E = "  +0.36%  " }
    @{ Row = 18; D = "2.648.66"; E = "  +1.12%  " }
    @{ Row = 19; D = "11.80"; E = "  +0.46%  " }
    @{ Row = 20; D = "7.84"; E = "  -0.05%  " }
    @{ Row = 21; D = "354.29"; E = "  -2.45%  " }
    @{ Row = 22; D = "4.28"; E = "  -2.34%  " }
    @{ Row = 23; D = "4.64"; E = "  -2.84%  " }
    @{ Row = 24; D = "10.56"; E = "  -2.96%  " }
    @{ Row = 25; D = $null; E = "  +0.15%  " }
    @{ Row = 26; D = "1.92"; E = "  -4.01%  " }
    @{ Row = 27; D = "69.02"; E = "  -2.11%  " }
    @{ Row = 28; D = "2.754.84"; E = "  -0.81%  " }
    @{ Row = 29; D = "1.00"; E = "  +0.14%  " }
    @{ Row = 30; D = "0.0₃0998"; E = "  -2.31%  " }
    @{ Row = 31; D = "540.80"; E = "  -2.37%  " }
    @{ Row = 32; D = "7.90"; E = "  -1.25%  " }
    @{ Row = 33; D = "1.34"; E = "  -3.28%  " }
    @{ Row = 34; D = $null; E = "  -2.06%  " }
    @{ Row = 35; D = $null; E = "  +2.38%  " }
    @{ Row = 36; D = $null; E = "  -0.01%  " }
    @{ Row = 37; D = "1.49"; E = "  -3.03%  " }
    @{ Row = 38; D = "157.58"; E = "  -0.09%  " }
    @{ Row = 39; D = "18.92"; E = "  -2.22%  " }
    @{ Row = 40; D = "0.364"; E = "  -1.84%  " }
    @{ Row = 41; D = "18.22"; E = "  +1.69%  " }
    @{ Row = 42; D = "1.80"; E = "  -0.96%  " }
    @{ Row = 43; D = "5.14"; E = $null }
    @{ Row = 45; D = "2.41"; E = "  -4.18%  " }
    @{ Row = 46; D = "0.0₆0302"; E = "  -0.51%  " }
    @{ Row = 47; D = "151.40"; E = "  -1.48%  " }
    @{ Row = 48; D = "0.575"; E = "  -2.94%  " }
    @{ Row = 49; D = "3.77"; E = "  -2.53%  " }
    @{ Row = 50; D = "1.70"; E = "  -1.60%  " }
    @{ Row = 51; D = "0.0768"; E = "  -1.08%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Leading apostrophe forces text interpretation so values such as
        # "1.00" or "5.20" keep their exact digits instead of being parsed as numbers.
        $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
